$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.418179035186768
$ws.Range("B1").Value = 1.974778413772583
$ws.Range("C1").Value = 3.510717391967773
$ws.Range("D1").Value = 3.648233652114868
$ws.Range("E1").Value = 0.8450916409492493
